$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column K header, reusing the same header formatting as column J (s="1")
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "intervention_type"

# Fill in intervention_type values for each row
$ws.Range("K2").Value = "DRUG"
$ws.Range("K3").Value = "DRUG"
$ws.Range("K4").Value = "DEVICE"
$ws.Range("K5").Value = "DRUG"
$ws.Range("K6").Value = "PROCEDURE"
$ws.Range("K7").Value = "DEVICE"
$ws.Range("K8").Value = "DRUG"
$ws.Range("K9").Value = "DRUG"
